$d = $word.ActiveDocument

# Question 7's answer key ends with a list of single-letter "ListParagraph"
# paragraphs (C, A, D, B, D, A, C). The author fixed the last answer
# choice from "C" to "D" (commit: "Sua dap an cau hoi").
$lastIndex = $d.Paragraphs.Count
$last = $d.Paragraphs($lastIndex)
$r = $last.Range
$r.MoveEnd(1, -1)
$r.Text = "D"

# Word tracks the location of the most recent edit with the hidden
# "_GoBack" bookmark; moving it here (off the previously-edited "b. Tang
# hieu suat" paragraph) matches real Word's behaviour after this edit.
$last2 = $d.Paragraphs($d.Paragraphs.Count)
$r2 = $last2.Range
$r2.MoveEnd(1, -1)
$d.Bookmarks.Add("_GoBack", $r2)
